$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.413.75"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "1.900.89"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Formula = "=""325.24"""
$ws.Range("E5").Value = "  -2.80%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Formula = "=""0.4809"""
$ws.Range("E7").Value = "  +2.95%  "

$ws.Range("D8").Formula = "=""0.4070"""

$ws.Range("D9").Formula = "=""0.08057"""
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").Formula = "=""1.003"""
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").Formula = "=""23.24"""
$ws.Range("E11").Value = "  +3.96%  "

$ws.Range("D12").Value = "1.866.98"
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").Formula = "=""5.943"""
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").Formula = "=""7.078"""
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").Formula = "=""89.76"""
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Formula = "=""0.06686"""
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Formula = "=""0.00001032"""

$ws.Range("D19").Formula = "=""17.60"""
$ws.Range("E19").Value = "  -1.48%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "29.443.35"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").Formula = "=""5.533"""
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").Formula = "=""11.78"""
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("D24").Formula = "=""2.161"""
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("D25").Value = "2.202.82"
$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("D26").Formula = "=""155.23"""
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").Formula = "=""19.80"""
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Formula = "=""6.062"""
$ws.Range("E28").Value = "  +5.46%  "

$ws.Range("D29").Formula = "=""2.091"""

$ws.Range("D30").Formula = "=""118.47"""
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").Formula = "=""1.030"""
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("D32").Formula = "=""0.09506"""
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Formula = "=""3.543"""
$ws.Range("E33").Value = "  -1.02%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Formula = "=""5.409"""
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Formula = "=""1.389"""
$ws.Range("E35").Value = "  -3.12%  "

$ws.Range("D36").Formula = "=""0.02250"""
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Formula = "=""0.06064"""

$ws.Range("D38").Formula = "=""1.176"""
$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Formula = "=""7.916"""
$ws.Range("E39").Value = "  -6.00%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Formula = "=""0.5867"""
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Formula = "=""0.1844"""
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Formula = "=""1.282"""
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Formula = "=""2.401"""
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("D45").Formula = "=""0.07755"""
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("D46").Formula = "=""12.31"""
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").Formula = "=""0.5520"""
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Formula = "=""1.919"""
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").Formula = "=""113.45"""
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("D50").Formula = "=""0.2943"""
$ws.Range("E50").Value = "  -1.91%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Formula = "=""43.67"""
$ws.Range("E51").Value = "  -1.12%  "

# Flatten any formula-based text assignments (used to avoid numeric auto-coercion)
# back into plain static text values, matching the source text cell type.
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$excel.CutCopyMode = 0
